# "Generate Report for Handback" - mark the localization status as handed
# back (in sync with en-US), refresh the handback timestamps, and clear the
# stale "handback file is not latest" error now that the report is current.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) ---
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("K2").Value2 = "2016-09-06 23:01:47"
$wsZhCn.Range("P2").Value2 = ""

# --- de-de sheet ---
$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("K2").Value2 = "2016-09-06 23:01:55"
$wsDeDe.Range("P2").Value2 = ""

# --- Column width adjustments to fit the new, longer status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
